$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "42×79=3318"
$t.Cell(1,2).Range.Text = "85×22=1870"
$t.Cell(1,3).Range.Text = "28×61=1708"
$t.Cell(1,4).Range.Text = "65×54=3510"
$t.Cell(1,5).Range.Text = "78×21=1638"

$t.Cell(5,1).Range.Text = "67×92=6164"
$t.Cell(5,2).Range.Text = "66×46=3036"
$t.Cell(5,3).Range.Text = "68×45=3060"
$t.Cell(5,4).Range.Text = "78×57=4446"
$t.Cell(5,5).Range.Text = "93×24=2232"

$t.Cell(10,1).Range.Text = "72×49=3528"
$t.Cell(10,2).Range.Text = "46×74=3404"
$t.Cell(10,3).Range.Text = "29×43=1247"
$t.Cell(10,4).Range.Text = "28×50=1400"
$t.Cell(10,5).Range.Text = "72×96=6912"

$t.Cell(15,1).Range.Text = "27×29=783"
$t.Cell(15,2).Range.Text = "27×18=486"
$t.Cell(15,3).Range.Text = "63×45=2835"
$t.Cell(15,4).Range.Text = "89×77=6853"
$t.Cell(15,5).Range.Text = "69×82=5658"

$t.Cell(20,1).Range.Text = "33×60=1980"
$t.Cell(20,2).Range.Text = "59×63=3717"
$t.Cell(20,3).Range.Text = "72×42=3024"
$t.Cell(20,4).Range.Text = "12×77=924"
$t.Cell(20,5).Range.Text = "64×32=2048"
